$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns keep their text formatting
# so numeric-looking strings (e.g. "21.48", "4.514") are not converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '26.214.89'
$ws.Range("E2").Value = '  -0.18%  '
$ws.Range("D3").Value = '1.682.79'
$ws.Range("E3").Value = '  +0.41%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '216.47'
$ws.Range("E5").Value = '  -0.35%  '
$ws.Range("D6").Value = '0.5238'
$ws.Range("E6").Value = '  -1.53%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '0.2694'
$ws.Range("E8").Value = '  +0.57%  '
$ws.Range("D9").Value = '0.06368'
$ws.Range("E9").Value = '  -1.61%  '
$ws.Range("D10").Value = '21.48'
$ws.Range("D11").Value = '0.07622'
$ws.Range("E11").Value = '  +1.42%  '
$ws.Range("D12").Value = '1.702.13'
$ws.Range("E12").Value = '  +1.43%  '
$ws.Range("D13").Value = '4.514'
$ws.Range("E13").Value = '  +0.09%  '
$ws.Range("D14").Value = '0.5754'
$ws.Range("E14").Value = '  -0.08%  '
$ws.Range("D15").Value = '0.000008331'
$ws.Range("E15").Value = '  -1.82%  '
$ws.Range("D16").Value = '66.04'
$ws.Range("E16").Value = '  +2.27%  '
$ws.Range("D17").Value = '26.257.77'
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("E18").Value = '  +0.01%  '
$ws.Range("D19").Value = '4.866'
$ws.Range("E19").Value = '  -0.84%  '
$ws.Range("E20").Value = '  -0.33%  '
$ws.Range("D21").Value = '189.69'
$ws.Range("E21").Value = '  -0.18%  '
$ws.Range("D22").Value = '6.234'
$ws.Range("E22").Value = '  +0.93%  '
$ws.Range("D24").Value = '148.89'
$ws.Range("E24").Value = '  +2.74%  '
$ws.Range("D25").Value = '7.787'
$ws.Range("E25").Value = '  -0.13%  '
$ws.Range("E26").Value = '  -1.03%  '
$ws.Range("D27").Value = '15.74'
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("D28").Value = '0.06284'
$ws.Range("E28").Value = '  -2.99%  '
$ws.Range("E29").Value = '  +0.92%  '
$ws.Range("D30").Value = '1.316'
$ws.Range("E30").Value = '  -0.09%  '
$ws.Range("D31").Value = '3.566'
$ws.Range("E31").Value = '  -0.40%  '
$ws.Range("D32").Value = '3.571'
$ws.Range("E32").Value = '  -0.33%  '
$ws.Range("D33").Value = '1.679'
$ws.Range("E33").Value = '  +1.62%  '
$ws.Range("D34").Value = '1.024'
$ws.Range("E34").Value = '  -0.50%  '
$ws.Range("D35").Value = '0.6120'
$ws.Range("E35").Value = '  -1.05%  '
$ws.Range("D36").Value = '2.419'
$ws.Range("E36").Value = '  +0.69%  '
$ws.Range("E37").Value = '  +1.13%  '
$ws.Range("D38").Value = '6.192'
$ws.Range("E38").Value = '  -1.58%  '
$ws.Range("E39").Value = '  -0.18%  '
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '1.097.81'
$ws.Range("E40").Value = '  -1.40%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '0.8872'
$ws.Range("E41").Value = '  +1.64%  '
$ws.Range("E42").Value = '  -0.35%  '
$ws.Range("D43").Value = '100.51'
$ws.Range("E43").Value = '  +0.07%  '
$ws.Range("D44").Value = '1.833.29'
$ws.Range("E44").Value = '  +0.39%  '
$ws.Range("D45").Value = '0.00000000110'
$ws.Range("E45").Value = '  +1.26%  '
$ws.Range("D46").Value = '57.31'
$ws.Range("E46").Value = '  +0.74%  '
$ws.Range("E47").Value = '  +0.54%  '
$ws.Range("D48").Value = '8.077'
$ws.Range("E48").Value = '  -1.18%  '
$ws.Range("D49").Value = '0.05278'
$ws.Range("E49").Value = '  +0.37%  '
$ws.Range("D50").Value = '0.4286'
$ws.Range("E50").Value = '  -0.06%  '
$ws.Range("D51").Value = '6.022'
$ws.Range("E51").Value = '  -0.84%  '
